$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand merged groups by inserting rows *inside* the existing merges ---
# Group 1 currently spans rows 2-4 (3 rows); needs to become rows 2-6 (5 rows): insert 2 rows.
$ws.Rows("4:5").Insert()
# Group 2 is now at rows 7-8 (2 rows); needs to become rows 7-10 (4 rows): insert 2 rows.
$ws.Rows("8:9").Insert()

# --- Row 2 (group 1 header row) ---
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Department of systema and computer engineering - admin staff "
$ws.Range("C2").Value = "Labiche"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Labiche"
$ws.Range("G2").Value = "Yvan"
$ws.Range("H2").Value = 0

# --- Row 3 ---
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = "McConnell"
$ws.Range("G3").Value = "Jenna"
$ws.Range("H3").Value = 0

# --- Row 4 ---
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = "Poll"
$ws.Range("G4").Value = "Jennifer"
$ws.Range("H4").Value = 0

# --- Row 5 ---
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Warmington"
$ws.Range("G5").Value = "Saundra"
$ws.Range("H5").Value = 0

# --- Row 6 ---
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = "East"
$ws.Range("G6").Value = "Erica"
$ws.Range("H6").Value = 0

# --- Row 7 (group 2 header row) ---
$ws.Range("A7").Value = 20
$ws.Range("B7").Value = "Department of systema and computer engineering - Tech  staff "
$ws.Range("C7").Value = "Labiche"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Buburuz"
$ws.Range("G7").Value = "Jerry "
$ws.Range("H7").Value = 0

# --- Row 8 ---
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "Russ"
$ws.Range("G8").Value = "Daren "
$ws.Range("H8").Value = 0

# --- Row 9 ---
$ws.Range("E9").Value = 17
$ws.Range("F9").Value = "Chiv"
$ws.Range("G9").Value = "kong "
$ws.Range("H9").Value = 0

# --- Row 10 ---
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = "Singh"
$ws.Range("G10").Value = "Ishdeep"
$ws.Range("H10").Value = 0

# --- Selection / view bookkeeping (cosmetic, matches author's final cursor position) ---
$ws.Range("H19").Select()
